$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($row, $col, $val) {
    $cell = $ws.Cells.Item($row, $col)
    $origStyle = $cell.Style
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.Style = $origStyle
}

$updates = @(
    @{ Row = 2;  D = "43.753.68";  E = "  -0.33%  " },
    @{ Row = 3;  D = "2.285.64";   E = "  -0.60%  " },
    @{ Row = 4;  D = $null;        E = "  +0.34%  " },
    @{ Row = 5;  D = "115.51";     E = "  +15.51%  " },
    @{ Row = 6;  D = "269.22";     E = "  -0.53%  " },
    @{ Row = 7;  D = $null;        E = "  -0.07%  " },
    @{ Row = 8;  D = $null;        E = "  +0.08%  " },
    @{ Row = 9;  D = "0.618";      E = "  +1.65%  " },
    @{ Row = 10; D = "48.58";      E = "  +8.26%  " },
    @{ Row = 11; D = "0.0942";     E = "  +1.25%  " },
    @{ Row = 12; D = "8.98";       E = "  +13.35%  " },
    @{ Row = 13; D = $null;        E = "  +0.41%  " },
    @{ Row = 14; D = "15.80";      E = "  -0.37%  " },
    @{ Row = 15; D = "2.632.25";   E = "  -0.40%  " },
    @{ Row = 16; D = "0.874";      E = "  +1.98%  " },
    @{ Row = 17; D = "2.295.65";   E = "  +0.11%  " },
    @{ Row = 18; D = "43.619.29";  E = "  -0.48%  " },
    @{ Row = 19; D = $null;        E = "  -1.21%  " },
    @{ Row = 20; D = "6.98";       E = "  +11.87%  " },
    @{ Row = 21; D = "72.46";      E = "  +0.15%  " },
    @{ Row = 22; D = $null;        E = "  -1.54%  " },
    @{ Row = 23; D = "10.19";      E = "  +11.56%  " },
    @{ Row = 24; D = "233.12";     E = "  -0.11%  " },
    @{ Row = 25; D = $null;        E = "  +2.80%  " },
    @{ Row = 26; D = "0.999";      E = $null },
    @{ Row = 27; D = "11.68";      E = "  +4.12%  " },
    @{ Row = 28; D = $null;        E = "  +13.13%  " },
    @{ Row = 29; D = "42.05";      E = "  +9.64%  " },
    @{ Row = 30; D = $null;        E = "  -1.97%  " },
    @{ Row = 31; D = $null;        E = "  -1.59%  " },
    @{ Row = 32; D = "173.48";     E = "  -1.94%  " },
    @{ Row = 33; D = "0.0933";     E = "  +4.33%  " },
    @{ Row = 34; D = "21.54";      E = "  -1.42%  " },
    @{ Row = 35; D = "5.72";       E = "  +4.93%  " },
    @{ Row = 36; D = "0.127";      E = "  +0.03%  " },
    @{ Row = 37; D = "4.75";       E = "  +0.23%  " },
    @{ Row = 38; D = "0.0360";     E = "  +2.08%  " },
    @{ Row = 39; D = $null;        E = "  -0.51%  " },
    @{ Row = 40; D = "3.84";       E = "  +8.63%  " },
    @{ Row = 41; D = "14.54";      E = "  +18.91%  " },
    @{ Row = 42; D = "74.78";      E = "  +15.16%  " },
    @{ Row = 43; D = "2.43";       E = "  +3.68%  " },
    @{ Row = 44; D = "0.241";      E = "  +1.78%  " },
    @{ Row = 45; D = "6.37";       E = "  +21.64%  " },
    @{ Row = 46; D = $null;        E = "  +0.19%  " },
    @{ Row = 47; D = "1.39";       E = "  +0.90%  " },
    @{ Row = 48; D = "8.75";       E = "  -1.13%  " },
    @{ Row = 49; D = "102.67";     E = "  +4.11%  " },
    @{ Row = 50; D = $null;        E = "  +3.26%  " }
)

foreach ($u in $updates) {
    if ($null -ne $u.D) {
        Set-TextValue $u.Row 4 $u.D
    }
    if ($null -ne $u.E) {
        Set-TextValue $u.Row 5 $u.E
    }
}
